# "20h 27m before PA3PbIB O4KA by Polyakov"
#
# The underlying commit is mostly a "re-save on a different machine" --
# absolute path, revision GUID, window geometry, default-row-height
# metrics and the localized "Normal"/"Обычный" cell-style name all shift
# as a side effect of Excel (on a different install) touching the file,
# and aren't reachable through the Excel object model (they're host/UI
# fingerprints, not document content). The only real, user-driven edits
# are: the data correction in G6, and where the cursor was left
# (selection) when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data fix: row 5 (employee "Чихватова Алёна Алексеевна") salary corrected
# from 2550 to 2600.
$ws.Range("G6").Value = 2600

# Cursor/selection left on C8 instead of G2 when the file was saved.
[void]$ws.Range("C8").Select()
